$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 16:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1572114
$ws.Range("C4").Value = 1531
$ws.Range("E4").Value = 1117293

# Alemania (row 11)
$ws.Range("B11").Value = 177989
$ws.Range("C11").Value = 162
$ws.Range("E11").Value = 12880
$ws.Range("G11").Value = 16
$ws.Range("H11").Value = 8209

# India (row 14)
$ws.Range("B14").Value = 107819
$ws.Range("C14").Value = 1344
$ws.Range("D14").Value = 43070
$ws.Range("E14").Value = 61432
$ws.Range("G14").Value = 15
$ws.Range("H14").Value = 3317

# Finlandia (row 61)
$ws.Range("E61").Value = 1139
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 304

# Kenia (row 105)
$ws.Range("D105").Value = 366
$ws.Range("E105").Value = 613
